$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("img_1.jpg")
$ws1.Range("B5").Value = 109
$ws1.Range("C5").Value = 88
$ws1.Range("D5").Value = 52.30021667480469
$ws1.Range("E5").Value = 66.3616943359375
$ws1.Range("F5").Value = 0.3898895084857941

$ws1.Range("D6").Value = 56.42783355712891
$ws1.Range("E6").Value = 66.79560089111328
$ws1.Range("F6").Value = 0.7200009822845459

$ws2 = $wb.Worksheets.Item("img_2.jpg")
$ws2.Range("B3").Value = -87
$ws2.Range("D3").Value = 54.81571960449219
$ws2.Range("E3").Value = 53.42763900756836
$ws2.Range("F3").Value = 0.5847430229187012

$ws2.Range("B8").Value = 99
$ws2.Range("C8").Value = 7
$ws2.Range("D8").Value = 54.98033142089844
$ws2.Range("E8").Value = 57.45763778686523
$ws2.Range("F8").Value = 0.6259744763374329
